# added environment specific credentials fetch
#
# Adds a new worksheet "test_login_visual" (a copy of the first sheet's
# layout/data) as the last tab, makes it the active sheet/selection, and
# points the workbook at its new on-disk test-data folder.

$wb = $excel.ActiveWorkbook

# Source sheet to clone the visual layout (values + styles + row heights) from.
$srcSheet = $wb.Worksheets.Item(1)

# New sheet goes after the current last tab.
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$newSheet = $wb.Worksheets.Add($null, $lastSheet)
$newSheet.Name = "test_login_visual"

# Copy values + formatting (shared strings / cell styles) from sheet1.
$srcSheet.Range("A1:E3").Copy($newSheet.Range("A1"))

# Match the source sheet's (taller, wrapped-text) row heights.
$newSheet.Rows.Item(1).RowHeight = $srcSheet.Rows.Item(1).RowHeight
$newSheet.Rows.Item(2).RowHeight = $srcSheet.Rows.Item(2).RowHeight
$newSheet.Rows.Item(3).RowHeight = $srcSheet.Rows.Item(3).RowHeight

# New sheet becomes the active tab/selection.
$newSheet.Range("E6").Select()

# Record the new absolute folder this workbook's test data now lives in.
$wb.Path = "C:\Users\ASUS\correctTesDataFile\playWrightWebAutomation\testData\"
